$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -13.3417
$ws.Range("B12").Value = 4.789499999999999
$ws.Range("C23").Value = -12.22900000000001
$ws.Range("C28").Value = -13.4195
$ws.Range("B32").Value = 6.826199999999997
$ws.Range("C32").Value = -11.9336
$ws.Range("C34").Value = -11.72570000000001
$ws.Range("B36").Value = 9.33270000000001
$ws.Range("B38").Value = 5.300399999999996
$ws.Range("C42").Value = -12.2312
$ws.Range("B46").Value = 6.499400000000005
$ws.Range("B54").Value = 5.719300000000001
$ws.Range("C54").Value = -12.5396
$ws.Range("B55").Value = 5.430599999999998
$ws.Range("B67").Value = 5.343699999999995
$ws.Range("B69").Value = 5.256099999999996
$ws.Range("B72").Value = 5.169800000000007
$ws.Range("B91").Value = 5.074699999999998
$ws.Range("C97").Value = -12.1676
$ws.Range("B99").Value = 6.036499999999997
$ws.Range("C99").Value = -11.96600000000001
$ws.Range("C101").Value = -12.9108
$ws.Range("B104").Value = 9.438800000000004
